# Scheduled-runner Sheets update: refresh computed market-profit columns
# (currentAveragePrice*, LevePrice*, LeveProfit*) across several leve rows.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 11
$ws.Range("K11").Value = 340.27274
$ws.Range("H11").Value = 340.27274
$ws.Range("M11").Value = -200.27274
$ws.Range("I11").Value = 340.27274
# Row 106
$ws.Range("J106").Value = 125002500
$ws.Range("K106").Value = 37039176
$ws.Range("M106").Value = -37038545
$ws.Range("I106").Value = 37039176
$ws.Range("N106").Value = -125003762
$ws.Range("H106").Value = 78433680
$ws.Range("L106").Value = 125002500
# Row 132
$ws.Range("M132").Value = -2664.928400000001
$ws.Range("I132").Value = 1731.6428
$ws.Range("H132").Value = 15153527
$ws.Range("L132").Value = 125007510
$ws.Range("J132").Value = 41669170
$ws.Range("K132").Value = 5194.928400000001
$ws.Range("N132").Value = -125012570
# Row 137
$ws.Range("J137").Value = 1974.8
$ws.Range("L137").Value = 5924.4
$ws.Range("K137").Value = 3726.9642
$ws.Range("N137").Value = -11024.4
$ws.Range("H137").Value = 1497.8372
$ws.Range("M137").Value = -1176.9642
$ws.Range("I137").Value = 1242.3214

$ws = $wb.Worksheets.Item("ARM")
# Row 4
$ws.Range("J4").Value = 1000
$ws.Range("L4").Value = 1000
$ws.Range("K4").Value = 0
$ws.Range("N4").Value = -1232
$ws.Range("H4").Value = 1000
$ws.Range("M4").ClearContents()
$ws.Range("I4").Value = 0
# Row 32
$ws.Range("J32").Value = 10199.889
$ws.Range("K32").Value = 4768.393
$ws.Range("N32").Value = -10773.889
$ws.Range("H32").Value = 6089.5674
$ws.Range("M32").Value = -4481.393
$ws.Range("I32").Value = 4768.393
$ws.Range("L32").Value = 10199.889
# Row 43
$ws.Range("J43").Value = 11305.167
$ws.Range("N43").Value = -11931.167
$ws.Range("H43").Value = 11305.167
$ws.Range("L43").Value = 11305.167
# Row 61
$ws.Range("K61").Value = 6003.522
$ws.Range("H61").Value = 6003.522
$ws.Range("M61").Value = -5791.522
$ws.Range("I61").Value = 6003.522
# Row 102
$ws.Range("K102").Value = 2059471.6
$ws.Range("M102").Value = -2057849.6
$ws.Range("I102").Value = 2059471.6
$ws.Range("H102").Value = 1951236.2
# Row 105
$ws.Range("N105").ClearContents()
$ws.Range("L105").Value = 0
$ws.Range("H105").Value = 0
$ws.Range("J105").Value = 0
# Row 132
$ws.Range("M132").Value = -5982.8948
$ws.Range("I132").Value = 2837.6316
$ws.Range("H132").Value = 3493.5312
$ws.Range("L132").Value = 13356.462
$ws.Range("J132").Value = 4452.154
$ws.Range("K132").Value = 8512.8948
$ws.Range("N132").Value = -18416.462
# Row 136
$ws.Range("K136").Value = 18010.566
$ws.Range("M136").Value = -15460.566
$ws.Range("I136").Value = 6003.522
$ws.Range("H136").Value = 6003.522

$ws = $wb.Worksheets.Item("BSM")
# Row 10
$ws.Range("M10").Value = -660
$ws.Range("I10").Value = 800
$ws.Range("H10").Value = 2900
$ws.Range("K10").Value = 800
# Row 18
$ws.Range("J18").Value = 7210
$ws.Range("L18").Value = 7210
$ws.Range("N18").Value = -8268
$ws.Range("H18").Value = 7210
# Row 94
$ws.Range("I94").Value = 1100.5385
$ws.Range("H94").Value = 1595.1052
$ws.Range("K94").Value = 1100.5385
$ws.Range("M94").Value = -649.5385000000001
# Row 134
$ws.Range("I134").Value = 4379.676
$ws.Range("H134").Value = 3739.9622
$ws.Range("L134").Value = 6781.875
$ws.Range("J134").Value = 2260.625
$ws.Range("K134").Value = 13139.028
$ws.Range("N134").Value = -11851.875
$ws.Range("M134").Value = -10604.028
# Row 138
$ws.Range("H138").Value = 0
$ws.Range("J138").Value = 0
$ws.Range("L138").Value = 0
$ws.Range("N138").ClearContents()
# Row 139
$ws.Range("K139").Value = 28709
$ws.Range("N139").Value = -60105.715
$ws.Range("J139").Value = 49825.715
$ws.Range("H139").Value = 47186.125
$ws.Range("M139").Value = -23569
$ws.Range("I139").Value = 28709
$ws.Range("L139").Value = 49825.715

$ws = $wb.Worksheets.Item("CRP")
# Row 22
$ws.Range("M22").Value = 35
$ws.Range("I22").Value = 315
$ws.Range("N22").Value = -984.6
$ws.Range("H22").Value = 304.14285
$ws.Range("J22").Value = 284.6
$ws.Range("L22").Value = 284.6
$ws.Range("K22").Value = 315
# Row 58
$ws.Range("L58").Value = 2268.611
$ws.Range("J58").Value = 2268.611
$ws.Range("K58").Value = 3876551.2
$ws.Range("N58").Value = -2674.611
$ws.Range("H58").Value = 2733320.2
$ws.Range("M58").Value = -3876348.2
$ws.Range("I58").Value = 3876551.2
# Row 94
$ws.Range("I94").Value = 4716.2
$ws.Range("H94").Value = 3729
$ws.Range("L94").Value = 3180.5557
$ws.Range("J94").Value = 3180.5557
$ws.Range("K94").Value = 4716.2
$ws.Range("N94").Value = -4082.5557
$ws.Range("M94").Value = -4265.2
# Row 132
$ws.Range("M132").Value = -19357724
$ws.Range("I132").Value = 6453418
$ws.Range("H132").Value = 5002101
$ws.Range("L132").Value = 9362.667000000001
$ws.Range("J132").Value = 3120.889
$ws.Range("K132").Value = 19360254
$ws.Range("N132").Value = -14422.667
# Row 134
$ws.Range("I134").Value = 14495870
$ws.Range("H134").Value = 8132401.5
$ws.Range("L134").Value = 3910.3335
$ws.Range("J134").Value = 1303.4445
$ws.Range("K134").Value = 43487610
$ws.Range("N134").Value = -8980.333500000001
$ws.Range("M134").Value = -43485075
# Row 136
$ws.Range("K136").Value = 11629653.6
$ws.Range("N136").Value = -11905.833
$ws.Range("M136").Value = -11627103.6
$ws.Range("I136").Value = 3876551.2
$ws.Range("L136").Value = 6805.833
$ws.Range("H136").Value = 2733320.2
$ws.Range("J136").Value = 2268.611
# Row 141
$ws.Range("N141").Value = -42742
$ws.Range("H141").Value = 32382
$ws.Range("J141").Value = 32382
$ws.Range("L141").Value = 32382

$ws = $wb.Worksheets.Item("CUL")
# Row 5
$ws.Range("M5").Value = -677.375
$ws.Range("I5").Value = 263.125
$ws.Range("N5").Value = -18174449
$ws.Range("H5").Value = 2423387.8
$ws.Range("L5").Value = 18174225
$ws.Range("J5").Value = 6058075
$ws.Range("K5").Value = 789.375
# Row 23
$ws.Range("M23").Value = -33333344
$ws.Range("I23").Value = 11111193
$ws.Range("N23").Value = -713.333345
$ws.Range("H23").Value = 5555637
$ws.Range("J23").Value = 81.111115
$ws.Range("L23").Value = 243.333345
$ws.Range("K23").Value = 33333579
# Row 75
$ws.Range("K75").Value = 937.9999799999999
$ws.Range("H75").Value = 20409124
$ws.Range("M75").Value = 60.00002000000006
$ws.Range("I75").Value = 312.66666
# Row 78
$ws.Range("K78").Value = 2813.99994
$ws.Range("H78").Value = 20409124
$ws.Range("M78").Value = 2178.00006
$ws.Range("I78").Value = 312.66666
# Row 107
$ws.Range("J107").Value = 649.3333
$ws.Range("L107").Value = 1947.9999
$ws.Range("K107").Value = 1114.00002
$ws.Range("N107").Value = -5787.9999
$ws.Range("M107").Value = 805.9999800000001
$ws.Range("I107").Value = 371.33334
$ws.Range("H107").Value = 569.9048
# Row 113
$ws.Range("K113").Value = 13638207
$ws.Range("J113").Value = 1111612.6
$ws.Range("N113").Value = -3339177.8
$ws.Range("M113").Value = -13636037
$ws.Range("I113").Value = 4546069
$ws.Range("H113").Value = 3000563.8
$ws.Range("L113").Value = 3334837.8
# Row 117
$ws.Range("H117").Value = 17555066
$ws.Range("J117").Value = 25649556
$ws.Range("L117").Value = 76948668
$ws.Range("N117").Value = -76955552
# Row 135
$ws.Range("I135").Value = 263.125
$ws.Range("N135").Value = -54527745
$ws.Range("H135").Value = 2423387.8
$ws.Range("M135").Value = 166.875
$ws.Range("L135").Value = 54522675
$ws.Range("J135").Value = 6058075
$ws.Range("K135").Value = 2368.125

$ws = $wb.Worksheets.Item("GSM")
# Row 97
$ws.Range("M97").Value = -78.88890000000004
$ws.Range("I97").Value = 574.8889
$ws.Range("H97").Value = 698.9286
$ws.Range("K97").Value = 574.8889
# Row 132
$ws.Range("M132").Value = -18521778.5
$ws.Range("I132").Value = 6174769.5
$ws.Range("H132").Value = 5954285
$ws.Range("L132").Value = 3600
$ws.Range("J132").Value = 1200
$ws.Range("K132").Value = 18524308.5
$ws.Range("N132").Value = -8660

$ws = $wb.Worksheets.Item("LTW")
# Row 22
$ws.Range("M22").Value = -17857141
$ws.Range("I22").Value = 17857436
$ws.Range("N22").Value = -2310.4615
$ws.Range("H22").Value = 3789296.5
$ws.Range("J22").Value = 1720.4615
$ws.Range("L22").Value = 1720.4615
$ws.Range("K22").Value = 17857436
# Row 27
$ws.Range("K27").Value = 17857436
$ws.Range("J27").Value = 1720.4615
$ws.Range("N27").Value = -1934.4615
$ws.Range("M27").Value = -17857329
$ws.Range("I27").Value = 17857436
$ws.Range("H27").Value = 3789296.5
$ws.Range("L27").Value = 1720.4615
# Row 61
$ws.Range("K61").Value = 2875
$ws.Range("H61").Value = 2875
$ws.Range("M61").Value = -2673
$ws.Range("I61").Value = 2875
# Row 113
$ws.Range("K113").Value = 2875
$ws.Range("M113").Value = -705
$ws.Range("I113").Value = 2875
$ws.Range("H113").Value = 2875
# Row 132
$ws.Range("M132").Value = -55185836
$ws.Range("I132").Value = 18396122
$ws.Range("H132").Value = 14819590
$ws.Range("L132").Value = 7585.7145
$ws.Range("J132").Value = 2528.5715
$ws.Range("K132").Value = 55188366
$ws.Range("N132").Value = -12645.7145

$ws = $wb.Worksheets.Item("WVR")
# Row 62
$ws.Range("K62").Value = 2900
$ws.Range("N62").ClearContents()
$ws.Range("M62").Value = -2276
$ws.Range("I62").Value = 2900
$ws.Range("H62").Value = 2900
$ws.Range("J62").Value = 0
$ws.Range("L62").Value = 0
# Row 65
$ws.Range("K65").Value = 14500
$ws.Range("N65").ClearContents()
$ws.Range("H65").Value = 2900
$ws.Range("M65").Value = -11380
$ws.Range("I65").Value = 2900
$ws.Range("L65").Value = 0
$ws.Range("J65").Value = 0
# Row 132
$ws.Range("M132").Value = -173.5769
$ws.Range("I132").Value = 901.1923
$ws.Range("H132").Value = 1163.7435
$ws.Range("L132").Value = 5066.5386
$ws.Range("J132").Value = 1688.8462
$ws.Range("K132").Value = 2703.5769
$ws.Range("N132").Value = -10126.5386
